$d = $word.ActiveDocument

$pairs = @(
    @("N = 90,222", "N = 90,237"),
    @("55.0 (20.0, 100.0)", "385.0 (140.0, 700.0)"),
    @("87,543 (97)", "87,558 (97)"),
    @("51,718 (57)", "51,723 (57)"),
    @("38,504 (43)", "38,514 (43)"),
    @("7,114 (7.9)", "7,118 (7.9)"),
    @("22,110 (25)", "22,114 (25)"),
    @("21,314 (24)", "21,318 (24)"),
    @("39,684 (44)", "39,687 (44)"),
    @("11,581 (13)", "11,584 (13)"),
    @("19,656 (22)", "19,661 (22)"),
    @("23,661 (26)", "23,663 (26)"),
    @("20,893 (23)", "20,894 (23)"),
    @("6,071 (6.7)", "6,072 (6.7)"),
    @("8,360 (9.3)", "8,363 (9.3)"),
    @("52,089 (58)", "52,094 (58)"),
    @("32,056 (36)", "32,065 (36)"),
    @("6,077 (6.7)", "6,078 (6.7)"),
    @("18,236 (20)", "18,243 (20)"),
    @("22,695 (25)", "22,696 (25)"),
    @("23,626 (26)", "23,629 (26)"),
    @("20,706 (23)", "20,710 (23)"),
    @("64,254 (71)", "64,263 (71)"),
    @("22,773 (25)", "22,778 (25)"),
    @("3,195 (3.5)", "3,196 (3.5)"),
    @("16,289 (18)", "16,292 (18)"),
    @("30,593 (34)", "30,597 (34)"),
    @("43,340 (48)", "43,348 (48)"),
    @("74,514 (83)", "74,526 (83)"),
    @("14,987 (17)", "14,990 (17)"),
    @("76,521 (85)", "76,532 (85)"),
    @("13,277 (15)", "13,281 (15)"),
    @("29,453 (33)", "29,456 (33)"),
    @("54,031 (60)", "54,040 (60)"),
    @("6,738 (7.5)", "6,741 (7.5)"),
    @("15,368 (17)", "15,369 (17)"),
    @("66,451 (74)", "66,462 (74)"),
    @("8,403 (9.3)", "8,406 (9.3)")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "WARNING: not found -> $old"
    }
}
